$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Week 5")

# Grades for week 5 (rows 2-12): Tutorial_05_1..6 (A2:A7), Program_05_1..5 (A8:A12)
# Column B = pass/fail flag (1 = passed), Column D = grader notes.
$ws.Range("B2:B12").Value = 1

# Stamp the (still blank) notes column with the sheet's normal/general format,
# matching the rest of the already-graded week sheets (D column cells carry an
# explicit "Normal" style even when empty).
$ws.Range("A1").Copy() | Out-Null
$ws.Range("D2:D12").PasteSpecial(-4122) | Out-Null

# Program_05_1 (row 8) failed to run.
$ws.Range("D8").Value = "Does not run"

# Narrow columns B and D to fit the newly entered data (they were sized for
# the still-empty template before grading).
$ws.Columns.Item(2).ColumnWidth = 8.333333333333334
$ws.Columns.Item(4).ColumnWidth = 12.333333333333334
